$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2976.25
$ws.Range("J112").Value = 1502
$ws.Range("L112").Value = 4506
$ws.Range("N112").Value = -6722
$ws.Range("H138").Value = 150762.14
$ws.Range("J138").Value = 5196.643
$ws.Range("L138").Value = 15589.929
$ws.Range("N138").Value = -25869.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20708.383
$ws.Range("I32").Value = 18652.867
$ws.Range("K32").Value = 18652.867
$ws.Range("M32").Value = -18365.867
$ws.Range("H74").Value = 4438.213
$ws.Range("I74").Value = 60747.5
$ws.Range("J74").Value = 1935.5778
$ws.Range("K74").Value = 60747.5
$ws.Range("L74").Value = 1935.5778
$ws.Range("M74").Value = -59873.5
$ws.Range("N74").Value = -3683.5778
$ws.Range("H77").Value = 4438.213
$ws.Range("I77").Value = 60747.5
$ws.Range("J77").Value = 1935.5778
$ws.Range("K77").Value = 303737.5
$ws.Range("L77").Value = 9677.888999999999
$ws.Range("M77").Value = -299369.5
$ws.Range("N77").Value = -18413.889
$ws.Range("H80").Value = 84999.71000000001
$ws.Range("J80").Value = 84999.71000000001
$ws.Range("L80").Value = 84999.71000000001
$ws.Range("N80").Value = -86995.71000000001
$ws.Range("H83").Value = 84999.71000000001
$ws.Range("J83").Value = 84999.71000000001
$ws.Range("L83").Value = 254999.13
$ws.Range("N83").Value = -264983.13
$ws.Range("H130").Value = 42150
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3278.2173
$ws.Range("I20").Value = 2728.8
$ws.Range("K20").Value = 2728.8
$ws.Range("M20").Value = -2481.8
$ws.Range("H82").Value = 24892.857
$ws.Range("J82").Value = 99990
$ws.Range("L82").Value = 99990
$ws.Range("N82").Value = -100756
$ws.Range("H85").Value = 24892.857
$ws.Range("J85").Value = 99990
$ws.Range("L85").Value = 99990
$ws.Range("N85").Value = -102642
$ws.Range("H86").Value = 6312
$ws.Range("I86").Value = 9240.200000000001
$ws.Range("J86").Value = 4059.5386
$ws.Range("K86").Value = 9240.200000000001
$ws.Range("L86").Value = 4059.5386
$ws.Range("M86").Value = -8117.200000000001
$ws.Range("N86").Value = -6305.5386
$ws.Range("H89").Value = 6312
$ws.Range("I89").Value = 9240.200000000001
$ws.Range("J89").Value = 4059.5386
$ws.Range("K89").Value = 46201
$ws.Range("L89").Value = 20297.693
$ws.Range("M89").Value = -40585
$ws.Range("N89").Value = -31529.693
$ws.Range("H94").Value = 1142.6072
$ws.Range("I94").Value = 1207.4615
$ws.Range("K94").Value = 1207.4615
$ws.Range("M94").Value = -756.4614999999999
$ws.Range("H97").Value = 11464
$ws.Range("I97").Value = 6673.143
$ws.Range("K97").Value = 6673.143
$ws.Range("M97").Value = -5682.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3055.3057
$ws.Range("I31").Value = 1640.9131
$ws.Range("K31").Value = 1640.9131
$ws.Range("M31").Value = -1345.9131
$ws.Range("H34").Value = 3055.3057
$ws.Range("I34").Value = 1640.9131
$ws.Range("K34").Value = 1640.9131
$ws.Range("M34").Value = -1438.9131
$ws.Range("H58").Value = 3185.8096
$ws.Range("J58").Value = 4418.7144
$ws.Range("L58").Value = 4418.7144
$ws.Range("N58").Value = -4824.7144
$ws.Range("H69").Value = 14500
$ws.Range("I69").Value = 14500
$ws.Range("K69").Value = 14500
$ws.Range("M69").Value = -13751
$ws.Range("H72").Value = 14500
$ws.Range("I72").Value = 14500
$ws.Range("K72").Value = 43500
$ws.Range("M72").Value = -39756
$ws.Range("H103").Value = 18312.908
$ws.Range("I103").Value = 18312.908
$ws.Range("K103").Value = 18312.908
$ws.Range("M103").Value = -17140.908
$ws.Range("H107").Value = 6942.9414
$ws.Range("I107").Value = 7314.375
$ws.Range("K107").Value = 7314.375
$ws.Range("M107").Value = -5394.375
$ws.Range("H136").Value = 3185.8096
$ws.Range("J136").Value = 4418.7144
$ws.Range("L136").Value = 13256.1432
$ws.Range("N136").Value = -18356.1432
$ws.Range("H141").Value = 412527.47
$ws.Range("J141").Value = 545563.25
$ws.Range("L141").Value = 545563.25
$ws.Range("N141").Value = -555923.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43144120
$ws.Range("I4").Value = 40538264
$ws.Range("J4").Value = 52590350
$ws.Range("K4").Value = 121614792
$ws.Range("L4").Value = 157771050
$ws.Range("M4").Value = -121614680
$ws.Range("N4").Value = -157771274
$ws.Range("H55").Value = 12100.733
$ws.Range("J55").Value = 13411.963
$ws.Range("L55").Value = 40235.889
$ws.Range("N55").Value = -40589.889
$ws.Range("H107").Value = 1223.1724
$ws.Range("I107").Value = 725.6667
$ws.Range("K107").Value = 2177.0001
$ws.Range("M107").Value = -257.0001000000002
$ws.Range("H113").Value = 2228.9333
$ws.Range("J113").Value = 2470.2222
$ws.Range("L113").Value = 7410.6666
$ws.Range("N113").Value = -11750.6666
$ws.Range("H139").Value = 1742096
$ws.Range("I139").Value = 2502181.8
$ws.Range("K139").Value = 7506545.399999999
$ws.Range("M139").Value = -7501405.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15326.357
$ws.Range("I80").Value = 18908.777
$ws.Range("K80").Value = 18908.777
$ws.Range("M80").Value = -17910.777
$ws.Range("H83").Value = 15326.357
$ws.Range("I83").Value = 18908.777
$ws.Range("K83").Value = 94543.88499999999
$ws.Range("M83").Value = -89551.88499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10028571
$ws.Range("I2").Value = 5599999
$ws.Range("K2").Value = 5599999
$ws.Range("M2").Value = -5599887
$ws.Range("H7").Value = 17519.89
$ws.Range("I7").Value = 21874.906
$ws.Range("J7").Value = 7565.5713
$ws.Range("K7").Value = 21874.906
$ws.Range("L7").Value = 7565.5713
$ws.Range("M7").Value = -21762.906
$ws.Range("N7").Value = -7789.5713
$ws.Range("H22").Value = 745.1905
$ws.Range("J22").Value = 768.5714
$ws.Range("L22").Value = 768.5714
$ws.Range("N22").Value = -1358.5714
$ws.Range("H27").Value = 745.1905
$ws.Range("J27").Value = 768.5714
$ws.Range("L27").Value = 768.5714
$ws.Range("N27").Value = -982.5714
$ws.Range("H61").Value = 24679
$ws.Range("I61").Value = 2856.4285
$ws.Range("K61").Value = 2856.4285
$ws.Range("M61").Value = -2654.4285
$ws.Range("H68").Value = 5000
$ws.Range("J68").Value = 5500
$ws.Range("L68").Value = 5500
$ws.Range("N68").Value = -6998
$ws.Range("H71").Value = 5000
$ws.Range("J71").Value = 5500
$ws.Range("L71").Value = 27500
$ws.Range("N71").Value = -34988
$ws.Range("H113").Value = 24679
$ws.Range("I113").Value = 2856.4285
$ws.Range("K113").Value = 2856.4285
$ws.Range("M113").Value = -686.4285
$ws.Range("H126").Value = 17519.89
$ws.Range("I126").Value = 21874.906
$ws.Range("J126").Value = 7565.5713
$ws.Range("K126").Value = 65624.71799999999
$ws.Range("L126").Value = 22696.7139
$ws.Range("M126").Value = -63154.71799999999
$ws.Range("N126").Value = -27636.7139
$ws.Range("H132").Value = 368208.06
$ws.Range("I132").Value = 714707.7
$ws.Range("K132").Value = 2144123.1
$ws.Range("M132").Value = -2141593.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 148129.86
$ws.Range("I62").Value = 232578.47
$ws.Range("J62").Value = 17618.363
$ws.Range("K62").Value = 232578.47
$ws.Range("L62").Value = 17618.363
$ws.Range("M62").Value = -231954.47
$ws.Range("N62").Value = -18866.363
$ws.Range("H65").Value = 148129.86
$ws.Range("I65").Value = 232578.47
$ws.Range("J65").Value = 17618.363
$ws.Range("K65").Value = 1162892.35
$ws.Range("L65").Value = 88091.815
$ws.Range("M65").Value = -1159772.35
$ws.Range("N65").Value = -94331.815
$ws.Range("H132").Value = 13323.655
$ws.Range("I132").Value = 16267.255
$ws.Range("K132").Value = 48801.765
$ws.Range("M132").Value = -46271.765
